$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23
$ws.Cells.Item($row, 1).Value = "2025-05-23 09:04:01"
$ws.Cells.Item($row, 2).Value = "Parclose"
$ws.Cells.Item($row, 3).Value = "Sortie - Demande 20250523_084920"
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 7
$ws.Cells.Item($row, 6).Value = 6
